$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Thursday row: DDB DCAN Tut -> DDB DCAN TUT (casing fix) ---
$ws.Range("C9").Value = "DDB DCAN TUT"

# --- Links Table: each subject gets its own distinct Google Meet link ---
$ws.Range("C14").Value = "https://meet.google.com/qnx-rvwj-sgx?pli=1&authuser=1"
$ws.Range("C15").Value = "https://meet.google.com/hwq-mqef-ewm?pli=1&authuser=1"
$ws.Range("C16").Value = "https://meet.google.com/ozp-cjab-vwg?pli=1&authuser=1"
$ws.Range("C17").Value = "https://meet.google.com/tmc-imzs-rag?pli=1&authuser=1"
$ws.Range("C18").Value = "https://meet.google.com/hwq-mqef-ewm?pli=1&authuser=1"
$ws.Range("C19").Value = "https://meet.google.com/zry-qpco-gbr?pli=1&authuser=1"
$ws.Range("C20").Value = "https://meet.google.com/ycb-jmtr-nve?pli=1&authuser=1"
$ws.Range("C21").Value = "https://meet.google.com/ufn-uuaf-drt?pli=1&authuser=1"
$ws.Range("C22").Value = "https://meet.google.com/ykr-oigv-rwd?pli=1&authuser=1"
$ws.Range("C23").Value = "https://meet.google.com/dfm-cbgj-pzs?pli=1&authuser=1"
$ws.Range("C24").Value = "https://meet.google.com/qnx-rvwj-sgx?pli=1&authuser=1"
$ws.Range("C25").Value = "https://meet.google.com/gfj-wxiq-uiz?pli=1&authuser=1"
$ws.Range("C26").Value = "https://timetablesysem2.netlify.app/#?pli=1&authuser=1"

# --- C27: first lecture link now a *working* (real) hyperlink ---
$ws.Range("C27").Value = "https://meet.google.com/new2"
$ws.Hyperlinks.Add($ws.Range("C27"), "https://meet.google.com/new2", "", "", "https://meet.google.com/new2") | Out-Null

# --- widen column C so the longer links are readable ---
$ws.Columns.Item(3).ColumnWidth = 51.6640625

# --- move selection to E10 (as left by the author) ---
$ws.Range("E10").Select() | Out-Null
